# Update the cryptocurrency price/volume table to the latest scraped snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a "Price" column cell while keeping it text
# (these values like "1.003" or "28.551.12" are display strings, not numbers,
# so we force Text format before/while writing, then drop back to the default
# "Normal" style so no extra formatting is left behind on the cell).
function Set-PriceCell($addr, $text) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Row 2
Set-PriceCell "D2" '28.551.12'
$ws.Range("E2").Value = '  +0.34%  '

# Row 3
Set-PriceCell "D3" '1.827.46'
$ws.Range("E3").Value = '  +0.01%  '

# Row 4
Set-PriceCell "D4" '1.003'
$ws.Range("E4").Value = '  +0.14%  '

# Row 5
Set-PriceCell "D5" '317.27'
$ws.Range("E5").Value = '  +0.48%  '

# Row 6
$ws.Range("E6").Value = '  +0.18%  '

# Row 7
Set-PriceCell "D7" '0.5165'
$ws.Range("E7").Value = '  +0.00%  '

# Row 8
Set-PriceCell "D8" '0.3888'
$ws.Range("E8").Value = '  -1.22%  '

# Row 9
Set-PriceCell "D9" '0.08423'
$ws.Range("E9").Value = '  +9.11%  '

# Row 10
$ws.Range("E10").Value = '  +0.79%  '

# Row 11
Set-PriceCell "D11" '41.97'

# Row 12
Set-PriceCell "D12" '6.425'
$ws.Range("E12").Value = '  +2.28%  '

# Row 13
Set-PriceCell "D13" '21.33'
$ws.Range("E13").Value = '  +1.26%  '

# Row 14
$ws.Range("E14").Value = '  +0.13%  '

# Row 15
Set-PriceCell "D15" '7.535'
$ws.Range("E15").Value = '  -0.53%  '

# Row 16
Set-PriceCell "D16" '1.827.36'
$ws.Range("E16").Value = '  +0.10%  '

# Row 17
Set-PriceCell "D17" '94.51'
$ws.Range("E17").Value = '  +1.37%  '

# Row 18
Set-PriceCell "D18" '0.00001135'
$ws.Range("E18").Value = '  +5.01%  '

# Row 19
Set-PriceCell "D19" '0.06630'
$ws.Range("E19").Value = '  -0.03%  '

# Row 20
Set-PriceCell "D20" '17.79'
$ws.Range("E20").Value = '  +0.48%  '

# Row 21
$ws.Range("E21").Value = '  +0.16%  '

# Row 22
$ws.Range("E22").Value = '  +0.26%  '

# Row 23
Set-PriceCell "D23" '28.590.46'
$ws.Range("E23").Value = '  +0.46%  '

# Row 24
Set-PriceCell "D24" '11.45'
$ws.Range("E24").Value = '  +2.80%  '

# Row 25
Set-PriceCell "D25" '2.280'
$ws.Range("E25").Value = '  +0.92%  '

# Row 26
Set-PriceCell "D26" '21.24'
$ws.Range("E26").Value = '  +3.00%  '

# Row 27
Set-PriceCell "D27" '159.78'
$ws.Range("E27").Value = '  +1.58%  '

# Row 28
Set-PriceCell "D28" '2.036.00'
$ws.Range("E28").Value = '  +0.06%  '

# Row 29
Set-PriceCell "D29" '2.413'
$ws.Range("E29").Value = '  -1.59%  '

# Row 30
Set-PriceCell "D30" '126.02'
$ws.Range("E30").Value = '  +0.86%  '

# Row 31
Set-PriceCell "D31" '0.1100'
$ws.Range("E31").Value = '  -0.02%  '

# Row 32
Set-PriceCell "D32" '1.101'
$ws.Range("E32").Value = '  -2.48%  '

# Row 33
$ws.Range("E33").Value = '  +1.54%  '

# Row 34
Set-PriceCell "D34" '0.07597'
$ws.Range("E34").Value = '  +5.76%  '

# Row 35
Set-PriceCell "D35" '3.673'
$ws.Range("E35").Value = '  +0.18%  '

# Row 36
$ws.Range("B36").Value = 'VeChain'
$ws.Range("C36").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-PriceCell "D36" '0.02389'
$ws.Range("E36").Value = '  +2.86%  '

# Row 37
$ws.Range("B37").Value = 'Algorand'
$ws.Range("C37").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-PriceCell "D37" '0.2230'
$ws.Range("E37").Value = '  -0.19%  '

# Row 38
Set-PriceCell "D38" '5.268'
$ws.Range("E38").Value = '  +2.23%  '

# Row 39
Set-PriceCell "D39" '8.760'
$ws.Range("E39").Value = '  -2.30%  '

# Row 40
$ws.Range("E40").Value = '  +2.22%  '

# Row 41
Set-PriceCell "D41" '11.49'
$ws.Range("E41").Value = '  +1.98%  '

# Row 42
Set-PriceCell "D42" '1.192'
$ws.Range("E42").Value = '  -0.03%  '

# Row 43
Set-PriceCell "D43" '1.402'
$ws.Range("E43").Value = '  +0.58%  '

# Row 44
$ws.Range("E44").Value = '  +0.94%  '

# Row 45
Set-PriceCell "D45" '0.6036'
$ws.Range("E45").Value = '  +2.21%  '

# Row 46
Set-PriceCell "D46" '3.781'
$ws.Range("E46").Value = '  +2.03%  '

# Row 47
Set-PriceCell "D47" '127.66'
$ws.Range("E47").Value = '  +2.64%  '

# Row 48
Set-PriceCell "D48" '2.002'
$ws.Range("E48").Value = '  +1.08%  '

# Row 49
Set-PriceCell "D49" '1.206'
$ws.Range("E49").Value = '  +1.94%  '

# Row 50
$ws.Range("E50").Value = '  +0.95%  '

# Row 51
$ws.Range("E51").Value = '  +1.31%  '
